$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A1").Value = "änderung nummer 2"
$ws.Range("A2").Select()
